$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 3.820425
$ws.Cells.Item(2, 8).Value = 11.461275
$ws.Cells.Item(2, 9).Value = 0.02049663039797357
$ws.Cells.Item(2, 10).Value = 0.02049663039797357
$ws.Cells.Item(2, 13).Value = 0.1825283333333333
$ws.Cells.Item(2, 14).Value = 0.547585
$ws.Cells.Item(2, 15).Value = 0.001028331058213739
$ws.Cells.Item(2, 16).Value = 0.001028331058213739
$ws.Cells.Item(2, 17).Value = 0.6973358078749999
$ws.Cells.Item(2, 18).Value = 6.276022270875
$ws.Cells.Item(2, 19).Value = 0.00002107732162696404
$ws.Cells.Item(2, 20).Value = 0.00002107732162696405

$ws.Cells.Item(3, 7).Value = 3.820425
$ws.Cells.Item(3, 8).Value = 11.461275
$ws.Cells.Item(3, 9).Value = 0.02049663039797357
$ws.Cells.Item(3, 10).Value = 0.02049663039797357
$ws.Cells.Item(3, 15).Value = 0.0001759459539160193
$ws.Cells.Item(3, 16).Value = 0.0001759459539160193
$ws.Cells.Item(3, 17).Value = 0.119313146225
$ws.Cells.Item(3, 18).Value = 1.073818316025
$ws.Cells.Item(3, 19).Value = 0.000003606299187435537
$ws.Cells.Item(3, 20).Value = 0.000003606299187435537

$ws.Cells.Item(4, 7).Value = 3.820425
$ws.Cells.Item(4, 8).Value = 11.461275
$ws.Cells.Item(4, 9).Value = 0.02049663039797357
$ws.Cells.Item(4, 10).Value = 0.02049663039797357
$ws.Cells.Item(4, 13).Value = 103.239782
$ws.Cells.Item(4, 14).Value = 309.719346
$ws.Cells.Item(4, 15).Value = 0.5816339432625932
$ws.Cells.Item(4, 16).Value = 0.5816339432625932
$ws.Cells.Item(4, 17).Value = 394.41984414735
$ws.Cells.Item(4, 18).Value = 3549.77859732615
$ws.Cells.Item(4, 19).Value = 0.0119215359619693
$ws.Cells.Item(4, 20).Value = 0.0119215359619693

$ws.Cells.Item(5, 7).Value = 3.820425
$ws.Cells.Item(5, 8).Value = 11.461275
$ws.Cells.Item(5, 9).Value = 0.02049663039797357
$ws.Cells.Item(5, 10).Value = 0.02049663039797357
$ws.Cells.Item(5, 13).Value = 0.04852733333333333
$ws.Cells.Item(5, 14).Value = 0.145582
$ws.Cells.Item(5, 15).Value = 0.0002733940705404138
$ws.Cells.Item(5, 16).Value = 0.0002733940705404139
$ws.Cells.Item(5, 17).Value = 0.18539503745
$ws.Cells.Item(5, 18).Value = 1.66855533705
$ws.Cells.Item(5, 19).Value = 0.000005603657216864375
$ws.Cells.Item(5, 20).Value = 0.000005603657216864377

$ws.Cells.Item(6, 7).Value = 3.820425
$ws.Cells.Item(6, 8).Value = 11.461275
$ws.Cells.Item(6, 9).Value = 0.02049663039797357
$ws.Cells.Item(6, 10).Value = 0.02049663039797357
$ws.Cells.Item(6, 13).Value = 73.99751433333334
$ws.Cells.Item(6, 14).Value = 221.992543
$ws.Cells.Item(6, 15).Value = 0.4168883856547366
$ws.Cells.Item(6, 16).Value = 0.4168883856547366
$ws.Cells.Item(6, 17).Value = 282.701953696925
$ws.Cells.Item(6, 18).Value = 2544.317583272325
$ws.Cells.Item(6, 19).Value = 0.008544807157973002
$ws.Cells.Item(6, 20).Value = 0.008544807157973002

$ws.Cells.Item(7, 9).Value = 0.7542622677884155
$ws.Cells.Item(7, 10).Value = 0.7542622677884157
$ws.Cells.Item(7, 13).Value = 0.1825283333333333
$ws.Cells.Item(7, 14).Value = 0.547585
$ws.Cells.Item(7, 15).Value = 0.001028331058213739
$ws.Cells.Item(7, 16).Value = 0.001028331058213739
$ws.Cells.Item(7, 17).Value = 25.66149057895222
$ws.Cells.Item(7, 18).Value = 230.95341521057
$ws.Cells.Item(7, 19).Value = 0.0007756313160055557
$ws.Cells.Item(7, 20).Value = 0.000775631316005556

$ws.Cells.Item(8, 9).Value = 0.7542622677884155
$ws.Cells.Item(8, 10).Value = 0.7542622677884157
$ws.Cells.Item(8, 15).Value = 0.0001759459539160193
$ws.Cells.Item(8, 16).Value = 0.0001759459539160193
$ws.Cells.Item(8, 19).Value = 0.0001327093942088927
$ws.Cells.Item(8, 20).Value = 0.0001327093942088927

$ws.Cells.Item(9, 9).Value = 0.7542622677884155
$ws.Cells.Item(9, 10).Value = 0.7542622677884157
$ws.Cells.Item(9, 13).Value = 103.239782
$ws.Cells.Item(9, 14).Value = 309.719346
$ws.Cells.Item(9, 15).Value = 0.5816339432625932
$ws.Cells.Item(9, 16).Value = 0.5816339432625932
$ws.Cells.Item(9, 17).Value = 14514.38603960708
$ws.Cells.Item(9, 18).Value = 130629.4743564637
$ws.Cells.Item(9, 19).Value = 0.4387045370679622
$ws.Cells.Item(9, 20).Value = 0.4387045370679623

$ws.Cells.Item(10, 9).Value = 0.7542622677884155
$ws.Cells.Item(10, 10).Value = 0.7542622677884157
$ws.Cells.Item(10, 13).Value = 0.04852733333333333
$ws.Cells.Item(10, 14).Value = 0.145582
$ws.Cells.Item(10, 15).Value = 0.0002733940705404138
$ws.Cells.Item(10, 16).Value = 0.0002733940705404139
$ws.Cells.Item(10, 17).Value = 6.822413180538222
$ws.Cells.Item(10, 18).Value = 61.40171862484399
$ws.Cells.Item(10, 19).Value = 0.0002062108316457186
$ws.Cells.Item(10, 20).Value = 0.0002062108316457186

$ws.Cells.Item(11, 9).Value = 0.7542622677884155
$ws.Cells.Item(11, 10).Value = 0.7542622677884157
$ws.Cells.Item(11, 13).Value = 73.99751433333334
$ws.Cells.Item(11, 14).Value = 221.992543
$ws.Cells.Item(11, 15).Value = 0.4168883856547366
$ws.Cells.Item(11, 16).Value = 0.4168883856547366
$ws.Cells.Item(11, 17).Value = 10403.24251174182
$ws.Cells.Item(11, 18).Value = 93629.18260567641
$ws.Cells.Item(11, 19).Value = 0.3144431791785932
$ws.Cells.Item(11, 20).Value = 0.3144431791785933

$ws.Cells.Item(12, 7).Value = 30.51067
$ws.Cells.Item(12, 8).Value = 91.53201
$ws.Cells.Item(12, 9).Value = 0.1636901460399144
$ws.Cells.Item(12, 10).Value = 0.1636901460399144
$ws.Cells.Item(12, 13).Value = 0.1825283333333333
$ws.Cells.Item(12, 14).Value = 0.547585
$ws.Cells.Item(12, 15).Value = 0.001028331058213739
$ws.Cells.Item(12, 16).Value = 0.001028331058213739
$ws.Cells.Item(12, 17).Value = 5.569061743983333
$ws.Cells.Item(12, 18).Value = 50.12155569585
$ws.Cells.Item(12, 19).Value = 0.0001683276610963866
$ws.Cells.Item(12, 20).Value = 0.0001683276610963867

$ws.Cells.Item(13, 7).Value = 30.51067
$ws.Cells.Item(13, 8).Value = 91.53201
$ws.Cells.Item(13, 9).Value = 0.1636901460399144
$ws.Cells.Item(13, 10).Value = 0.1636901460399144
$ws.Cells.Item(13, 15).Value = 0.0001759459539160193
$ws.Cells.Item(13, 16).Value = 0.0001759459539160193
$ws.Cells.Item(13, 17).Value = 0.9528583943233333
$ws.Cells.Item(13, 18).Value = 8.57572554891
$ws.Cells.Item(13, 19).Value = 0.00002880061889164525
$ws.Cells.Item(13, 20).Value = 0.00002880061889164525

$ws.Cells.Item(14, 7).Value = 30.51067
$ws.Cells.Item(14, 8).Value = 91.53201
$ws.Cells.Item(14, 9).Value = 0.1636901460399144
$ws.Cells.Item(14, 10).Value = 0.1636901460399144
$ws.Cells.Item(14, 13).Value = 103.239782
$ws.Cells.Item(14, 14).Value = 309.719346
$ws.Cells.Item(14, 15).Value = 0.5816339432625932
$ws.Cells.Item(14, 16).Value = 0.5816339432625932
$ws.Cells.Item(14, 17).Value = 3149.91491947394
$ws.Cells.Item(14, 18).Value = 28349.23427526546
$ws.Cells.Item(14, 19).Value = 0.09520774511442519
$ws.Cells.Item(14, 20).Value = 0.09520774511442519

$ws.Cells.Item(15, 7).Value = 30.51067
$ws.Cells.Item(15, 8).Value = 91.53201
$ws.Cells.Item(15, 9).Value = 0.1636901460399144
$ws.Cells.Item(15, 10).Value = 0.1636901460399144
$ws.Cells.Item(15, 13).Value = 0.04852733333333333
$ws.Cells.Item(15, 14).Value = 0.145582
$ws.Cells.Item(15, 15).Value = 0.0002733940705404138
$ws.Cells.Item(15, 16).Value = 0.0002733940705404139
$ws.Cells.Item(15, 17).Value = 1.480601453313333
$ws.Cells.Item(15, 18).Value = 13.32541307982
$ws.Cells.Item(15, 19).Value = 0.00004475191533320701
$ws.Cells.Item(15, 20).Value = 0.00004475191533320702

$ws.Cells.Item(16, 7).Value = 30.51067
$ws.Cells.Item(16, 8).Value = 91.53201
$ws.Cells.Item(16, 9).Value = 0.1636901460399144
$ws.Cells.Item(16, 10).Value = 0.1636901460399144
$ws.Cells.Item(16, 13).Value = 73.99751433333334
$ws.Cells.Item(16, 14).Value = 221.992543
$ws.Cells.Item(16, 15).Value = 0.4168883856547366
$ws.Cells.Item(16, 16).Value = 0.4168883856547366
$ws.Cells.Item(16, 17).Value = 2257.713740644604
$ws.Cells.Item(16, 18).Value = 20319.42366580143
$ws.Cells.Item(16, 19).Value = 0.06824052073016801
$ws.Cells.Item(16, 20).Value = 0.06824052073016801

$ws.Cells.Item(17, 7).Value = 0.258813
$ws.Cells.Item(17, 8).Value = 0.776439
$ws.Cells.Item(17, 9).Value = 0.001388535150720334
$ws.Cells.Item(17, 10).Value = 0.001388535150720334
$ws.Cells.Item(17, 13).Value = 0.1825283333333333
$ws.Cells.Item(17, 14).Value = 0.547585
$ws.Cells.Item(17, 15).Value = 0.001028331058213739
$ws.Cells.Item(17, 16).Value = 0.001028331058213739
$ws.Cells.Item(17, 17).Value = 0.047240705535
$ws.Cells.Item(17, 18).Value = 0.425166349815
$ws.Cells.Item(17, 19).Value = 0.000001427873820907214
$ws.Cells.Item(17, 20).Value = 0.000001427873820907214

$ws.Cells.Item(18, 7).Value = 0.258813
$ws.Cells.Item(18, 8).Value = 0.776439
$ws.Cells.Item(18, 9).Value = 0.001388535150720334
$ws.Cells.Item(18, 10).Value = 0.001388535150720334
$ws.Cells.Item(18, 15).Value = 0.0001759459539160193
$ws.Cells.Item(18, 16).Value = 0.0001759459539160193
$ws.Cells.Item(18, 17).Value = 0.008082816261000001
$ws.Cells.Item(18, 18).Value = 0.072745346349
$ws.Cells.Item(18, 19).Value = 0.0000002443071416394128
$ws.Cells.Item(18, 20).Value = 0.0000002443071416394127

$ws.Cells.Item(19, 7).Value = 0.258813
$ws.Cells.Item(19, 8).Value = 0.776439
$ws.Cells.Item(19, 9).Value = 0.001388535150720334
$ws.Cells.Item(19, 10).Value = 0.001388535150720334
$ws.Cells.Item(19, 13).Value = 103.239782
$ws.Cells.Item(19, 14).Value = 309.719346
$ws.Cells.Item(19, 15).Value = 0.5816339432625932
$ws.Cells.Item(19, 16).Value = 0.5816339432625932
$ws.Cells.Item(19, 17).Value = 26.719797698766
$ws.Cells.Item(19, 18).Value = 240.478179288894
$ws.Cells.Item(19, 19).Value = 0.0008076191750721872
$ws.Cells.Item(19, 20).Value = 0.0008076191750721871

$ws.Cells.Item(20, 7).Value = 0.258813
$ws.Cells.Item(20, 8).Value = 0.776439
$ws.Cells.Item(20, 9).Value = 0.001388535150720334
$ws.Cells.Item(20, 10).Value = 0.001388535150720334
$ws.Cells.Item(20, 13).Value = 0.04852733333333333
$ws.Cells.Item(20, 14).Value = 0.145582
$ws.Cells.Item(20, 15).Value = 0.0002733940705404138
$ws.Cells.Item(20, 16).Value = 0.0002733940705404139
$ws.Cells.Item(20, 17).Value = 0.012559504722
$ws.Cells.Item(20, 18).Value = 0.113035542498
$ws.Cells.Item(20, 19).Value = 0.0000003796172769438792
$ws.Cells.Item(20, 20).Value = 0.0000003796172769438792

$ws.Cells.Item(21, 7).Value = 0.258813
$ws.Cells.Item(21, 8).Value = 0.776439
$ws.Cells.Item(21, 9).Value = 0.001388535150720334
$ws.Cells.Item(21, 10).Value = 0.001388535150720334
$ws.Cells.Item(21, 13).Value = 73.99751433333334
$ws.Cells.Item(21, 14).Value = 221.992543
$ws.Cells.Item(21, 15).Value = 0.4168883856547366
$ws.Cells.Item(21, 16).Value = 0.4168883856547366
$ws.Cells.Item(21, 17).Value = 19.151518677153
$ws.Cells.Item(21, 18).Value = 172.363668094377
$ws.Cells.Item(21, 19).Value = 0.0005788641774086566
$ws.Cells.Item(21, 20).Value = 0.0005788641774086565

$ws.Cells.Item(22, 7).Value = 11.213844
$ws.Cells.Item(22, 8).Value = 33.641532
$ws.Cells.Item(22, 9).Value = 0.0601624206229761
$ws.Cells.Item(22, 10).Value = 0.0601624206229761
$ws.Cells.Item(22, 13).Value = 0.1825283333333333
$ws.Cells.Item(22, 14).Value = 0.547585
$ws.Cells.Item(22, 15).Value = 0.001028331058213739
$ws.Cells.Item(22, 16).Value = 0.001028331058213739
$ws.Cells.Item(22, 17).Value = 2.04684425558
$ws.Cells.Item(22, 18).Value = 18.42159830022
$ws.Cells.Item(22, 19).Value = 0.00006186688566392507
$ws.Cells.Item(22, 20).Value = 0.00006186688566392507

$ws.Cells.Item(23, 7).Value = 11.213844
$ws.Cells.Item(23, 8).Value = 33.641532
$ws.Cells.Item(23, 9).Value = 0.0601624206229761
$ws.Cells.Item(23, 10).Value = 0.0601624206229761
$ws.Cells.Item(23, 15).Value = 0.0001759459539160193
$ws.Cells.Item(23, 16).Value = 0.0001759459539160193
$ws.Cells.Item(23, 17).Value = 0.350212086068
$ws.Cells.Item(23, 18).Value = 3.151908774612
$ws.Cells.Item(23, 19).Value = 0.00001058533448640632
$ws.Cells.Item(23, 20).Value = 0.00001058533448640632

$ws.Cells.Item(24, 7).Value = 11.213844
$ws.Cells.Item(24, 8).Value = 33.641532
$ws.Cells.Item(24, 9).Value = 0.0601624206229761
$ws.Cells.Item(24, 10).Value = 0.0601624206229761
$ws.Cells.Item(24, 13).Value = 103.239782
$ws.Cells.Item(24, 14).Value = 309.719346
$ws.Cells.Item(24, 15).Value = 0.5816339432625932
$ws.Cells.Item(24, 16).Value = 0.5816339432625932
$ws.Cells.Item(24, 17).Value = 1157.714809942008
$ws.Cells.Item(24, 18).Value = 10419.43328947807
$ws.Cells.Item(24, 19).Value = 0.03499250594316435
$ws.Cells.Item(24, 20).Value = 0.03499250594316435

$ws.Cells.Item(25, 7).Value = 11.213844
$ws.Cells.Item(25, 8).Value = 33.641532
$ws.Cells.Item(25, 9).Value = 0.0601624206229761
$ws.Cells.Item(25, 10).Value = 0.0601624206229761
$ws.Cells.Item(25, 13).Value = 0.04852733333333333
$ws.Cells.Item(25, 14).Value = 0.145582
$ws.Cells.Item(25, 15).Value = 0.0002733940705404138
$ws.Cells.Item(25, 16).Value = 0.0002733940705404139
$ws.Cells.Item(25, 17).Value = 0.544177945736
$ws.Cells.Item(25, 18).Value = 4.897601511623999
$ws.Cells.Item(25, 19).Value = 0.00001644804906767998
$ws.Cells.Item(25, 20).Value = 0.00001644804906767998

$ws.Cells.Item(26, 7).Value = 11.213844
$ws.Cells.Item(26, 8).Value = 33.641532
$ws.Cells.Item(26, 9).Value = 0.0601624206229761
$ws.Cells.Item(26, 10).Value = 0.0601624206229761
$ws.Cells.Item(26, 13).Value = 73.99751433333334
$ws.Cells.Item(26, 14).Value = 221.992543
$ws.Cells.Item(26, 15).Value = 0.4168883856547366
$ws.Cells.Item(26, 16).Value = 0.4168883856547366
$ws.Cells.Item(26, 17).Value = 829.7965821217641
$ws.Cells.Item(26, 18).Value = 7468.169239095876
$ws.Cells.Item(26, 19).Value = 0.02508101441059374
$ws.Cells.Item(26, 20).Value = 0.02508101441059374
